# Applies the edits described in the commit:
#   "Diagramas de robustez y secuencia CU 12 y 15"
# which updates the "Estado"/"Incremento" fields of CU-12 (CRU grupo) and
# CU-15 (Consultar grupos y rentas) to "planificado" / 1, and updates the
# description of CU-15 to mention "rentas" as well as "grupos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU-12 (row 16, "CRU grupo"): Estado vacio -> planificado, Incremento 0 -> 1
$ws.Range("E16").Value2 = "planificado"
$ws.Range("F16").Value2 = 1

# CU-15 (row 19, "Consultar grupos y rentas"):
#  - Descripcion updated to mention "grupos y rentas" instead of just "grupos"
#  - Estado vacio -> planificado, Incremento 0 -> 1
$ws.Range("C19").Value2 = "El director debera visualizar todos los grupos y rentas existentes en el sistema y los alumnos"
$ws.Range("E19").Value2 = "planificado"
$ws.Range("F19").Value2 = 1

# The longer description now wraps onto more lines, so the row grows taller
$ws.Rows.Item(19).RowHeight = 33.75

# Update the selected cell shown in the sheet view to E19
$ws.Range("E19").Select()
